# Amend corrected label annotations
# Normalize the case of various "labels" (column F) entries on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F13").Value  = "use restrictions"
$ws.Range("F20").Value  = "93_referral_statement"
$ws.Range("F26").Value  = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F28").Value  = "ppe"
$ws.Range("F29").Value  = "ppe"
$ws.Range("F31").Value  = "ppe"
$ws.Range("F36").Value  = "application instructions || env warning - species"
$ws.Range("F38").Value  = "env warning - water"
$ws.Range("F44").Value  = "application instructions"
$ws.Range("F46").Value  = "135_product_information"
$ws.Range("F48").Value  = "use restrictions"
$ws.Range("F49").Value  = "use restrictions"
$ws.Range("F50").Value  = "use restrictions"
$ws.Range("F67").Value  = "application instructions"
$ws.Range("F68").Value  = "application instructions"
$ws.Range("F69").Value  = "application instructions"
$ws.Range("F70").Value  = "use restrictions"
$ws.Range("F71").Value  = "application instructions"
$ws.Range("F72").Value  = "application instructions"
$ws.Range("F73").Value  = "application instructions"
$ws.Range("F74").Value  = "application instructions"
$ws.Range("F76").Value  = "use restrictions"
$ws.Range("F121").Value = "154_pesticide_storage"
